# Apply "cryptos list" update (GitHub Actions scheduled refresh).
# Updates Price (column D) and Volume(1h) (column E) values, and for the
# two pairs of rows whose rank swapped (28/29 and 34/35) also updates the
# Coin (B) and Link (C) columns.
#
# Column D holds prices that look numeric (e.g. "134.23") but must stay as
# plain text, exactly like the source workbook (every cell is an inline
# string, never a number - note values such as "64.293.40" that contain two
# '.' separators, which only make sense as text). Excel's COM Value setter
# auto-coerces numeric-looking strings into real numbers (losing trailing
# zeros, introducing float noise, etc.), so every column-D write below first
# forces the cell's number format to Text ("@") to keep it a literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "64.293.40"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.501.02"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "585.75"
$ws.Range("E5").Value = "  +0.44%  "

# Row 6 - Solana
Set-TextValue "D6" "134.23"
$ws.Range("E6").Value = "  +2.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.69%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.51%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.78%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.387"
$ws.Range("E11").Value = "  +1.56%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "4.099.68"
$ws.Range("E12").Value = "  +0.51%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.10%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +2.97%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.501.85"
$ws.Range("E15").Value = "  +0.82%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  -4.80%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "64.323.45"
$ws.Range("E17").Value = "  +0.88%  "

# Row 18 - Uniswap
Set-TextValue "D18" "9.94"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19 - Polkadot
Set-TextValue "D19" "5.75"
$ws.Range("E19").Value = "  +1.62%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.73"
$ws.Range("E20").Value = "  -3.73%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "395.21"
$ws.Range("E21").Value = "  +3.09%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -0.68%  "

# Row 23 - WrappedeETH
Set-TextValue "D23" "3.643.06"
$ws.Range("E23").Value = "  +0.62%  "

# Row 24 - Litecoin
Set-TextValue "D24" "74.18"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 - LEO
Set-TextValue "D26" "5.66"
$ws.Range("E26").Value = "  -0.69%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000115"
$ws.Range("E27").Value = "  +1.03%  "

# Row 28 - now Binance-PegBSC-USD (was RenderToken); ranking swapped with row 29
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - now RenderToken (was Binance-PegBSC-USD)
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "7.42"
$ws.Range("E29").Value = "  -1.52%  "

# Row 30 - Fetch.AI
Set-TextValue "D30" "1.50"
$ws.Range("E30").Value = "  -4.81%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 - PancakeSwap
Set-TextValue "D32" "2.23"
$ws.Range("E32").Value = "  +0.08%  "

# Row 33 - RenzoRestakedETH
Set-TextValue "D33" "3.522.22"
$ws.Range("E33").Value = "  +0.98%  "

# Row 34 - now Kaspa (was USDe); ranking swapped with row 35
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D34" "0.150"
$ws.Range("E34").Value = "  +4.18%  "

# Row 35 - now USDe (was Kaspa)
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "23.48"
$ws.Range("E36").Value = "  +0.17%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  -2.11%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +0.47%  "

# Row 39 - Aptos
Set-TextValue "D39" "6.91"
$ws.Range("E39").Value = "  -0.58%  "

# Row 40 - Monero
Set-TextValue "D40" "161.33"
$ws.Range("E40").Value = "  +1.28%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0784"
$ws.Range("E41").Value = "  -1.55%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.808"
$ws.Range("E42").Value = "  -0.04%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.06%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "25.26"
$ws.Range("E44").Value = "  -3.48%  "

# Row 45 - Filecoin
Set-TextValue "D45" "4.41"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46 - ONDO
Set-TextValue "D46" "1.18"
$ws.Range("E46").Value = "  -2.69%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  +2.03%  "

# Row 48 - Maker
Set-TextValue "D48" "2.472.52"
$ws.Range("E48").Value = "  +2.29%  "

# Row 49 - Cosmos
Set-TextValue "D49" "6.79"
$ws.Range("E49").Value = "  -0.74%  "

# Row 50 - SuiNetwork
$ws.Range("E50").Value = "  -0.72%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -1.05%  "
